$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# for the coin rows whose figures moved since the last scrape.
# A leading apostrophe is used for D-column values that look like plain
# numbers (single decimal point) so Excel keeps them as text, matching
# the original inlineStr/text cells instead of letting them be
# auto-converted to numeric values.

$ws.Range("D2").Value = '26.668.34'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '1.594.49'
$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("D5").Value = "'211.50"
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -1.88%  '
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("D10").Value = "'19.68"
$ws.Range("E10").Value = '  -1.63%  '
$ws.Range("D12").Value = '1.818.12'
$ws.Range("E12").Value = '  -1.82%  '
$ws.Range("D13").Value = '1.585.75'
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("E14").Value = '  -2.53%  '
$ws.Range("D16").Value = "'64.77"
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").Value = '26.644.16'
$ws.Range("E17").Value = '  -1.53%  '
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").Value = "'209.00"
$ws.Range("E19").Value = '  -1.89%  '
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("E21").Value = '  -2.23%  '
$ws.Range("D22").Value = "'4.24"
$ws.Range("E22").Value = '  -2.34%  '
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("D24").Value = "'8.87"
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("D25").Value = "'146.69"
$ws.Range("E25").Value = '  -0.89%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = "'7.13"
$ws.Range("E27").Value = '  -3.11%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").Value = "'15.33"
$ws.Range("E29").Value = '  -1.27%  '
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = '  -1.33%  '
$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = '  -3.68%  '
$ws.Range("D33").Value = "'0.661"
$ws.Range("E33").Value = '  -8.59%  '
$ws.Range("D35").Value = '1.287.33'
$ws.Range("E35").Value = '  -5.41%  '
$ws.Range("E36").Value = '  -0.77%  '
$ws.Range("E37").Value = '  -5.58%  '
$ws.Range("E38").Value = '  -3.08%  '
$ws.Range("D39").Value = "'0.833"
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("D41").Value = "'0.793"
$ws.Range("E41").Value = '  -1.18%  '
$ws.Range("E42").Value = '  +0.49%  '
$ws.Range("D43").Value = "'2.20"
$ws.Range("E43").Value = '  -1.51%  '
$ws.Range("D44").Value = "'63.58"
$ws.Range("E44").Value = '  -1.36%  '
$ws.Range("D45").Value = '1.730.32'
$ws.Range("E45").Value = '  -1.88%  '
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("D48").Value = "'0.870"
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("E50").Value = '  -1.84%  '
$ws.Range("D51").Value = "'7.48"
$ws.Range("E51").Value = '  -2.26%  '
